$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: trim the old "MYPROFILE" test case down to just the ID + tester name ---
# B2 keeps its style but loses its value (Car_SRS_45 is gone).
$ws.Range("B2").ClearContents()
# C2, D2, F2, G2, H2 held the old requirement/title/data/steps/ER text - remove them completely.
$ws.Range("C2").Clear()
$ws.Range("D2").Clear()
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()
# J2 now carries the "Designed by" name for the remaining test case.
$ws.Range("J2").Value = "Fatma"

# --- New rows 3-20: one new Car_Details test case id per row in column A ---
for ($i = 3; $i -le 20; $i++) {
  $n = $i - 1
  $id = "Car_Details_" + ("{0:D2}" -f $n)
  $ws.Cells.Item($i, 1).Value = $id
}

# Match the row height used throughout the sheet for the newly added rows.
$ws.Range("A3:A20").RowHeight = 89.25

# --- Remove the logo picture that used to sit to the right of the table ---
foreach ($shp in $ws.Shapes) {
  $shp.Delete()
}

# --- Update the view: scroll over to show column H / select H2:I2 ---
$ws.Range("H2:I2").Select()
$excel.ActiveWindow.ScrollColumn = 8
